$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 136
$ws1.Range("F7").Value = 13108
$ws1.Range("F10").Value = 274
$ws1.Range("F11").Value = 3513
$ws1.Range("F13").Value = 6615
$ws1.Range("F16").Value = 3487
$ws1.Range("F21").Value = 70
$ws1.Range("F22").Value = 126
$ws1.Range("F23").Value = 54
$ws1.Range("F24").Value = 3646
$ws1.Range("F27").Value = 3263
$ws1.Range("F29").Value = 1910
$ws1.Range("F31").Value = 233
$ws1.Range("F32").Value = 6796
$ws1.Range("F33").Value = 168
$ws1.Range("F34").Value = 1417
$ws1.Range("F35").Value = 2018
$ws1.Range("F36").Value = 1302
$ws1.Range("F37").Value = 109
$ws1.Range("F40").Value = 218
$ws1.Range("F44").Value = 144
$ws1.Range("F45").Value = 1222
$ws1.Range("F46").Value = 1808
$ws1.Range("F49").Value = 1179

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 927

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 454
$ws3.Range("F3").Value = 623
$ws3.Range("F4").Value = 26

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 136
$ws4.Range("F6").Value = 454
$ws4.Range("F7").Value = 623
$ws4.Range("F9").Value = 13108
$ws4.Range("F13").Value = 274
$ws4.Range("F14").Value = 3514
$ws4.Range("F16").Value = 3487
$ws4.Range("F20").Value = 70
$ws4.Range("F22").Value = 126
$ws4.Range("F23").Value = 54
$ws4.Range("F24").Value = 3646
$ws4.Range("F27").Value = 3264
$ws4.Range("F28").Value = 3264
$ws4.Range("F30").Value = 1910
$ws4.Range("F32").Value = 233
$ws4.Range("F33").Value = 6796
$ws4.Range("F35").Value = 168
$ws4.Range("F36").Value = 1417
$ws4.Range("F37").Value = 2018
$ws4.Range("F39").Value = 1302
$ws4.Range("F40").Value = 109
$ws4.Range("F42").Value = 218
$ws4.Range("F45").Value = 1222
$ws4.Range("F47").Value = 1808
